# Auto-applies the cryptos-list refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.820.20"
$ws.Range("E2").Value = "  +0.58%  "

# Row 3: Ethereum -> Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.809.16"
$ws.Range("E3").Value = "  +1.07%  "

# Row 4: TetherUSD -> TetherUSD
$ws.Range("E4").Value = "  -0.08%  "

# Row 5: BNB -> BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.45"
$ws.Range("E5").Value = "  +0.87%  "

# Row 6: Solana -> Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.62"
$ws.Range("E6").Value = "  +0.70%  "

# Row 7: LidoStakedEther -> LidoStakedEther
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.806.75"
$ws.Range("E7").Value = "  +1.04%  "

# Row 8: USDC -> USDC
$ws.Range("E8").Value = "  +0.01%  "

# Row 9: XRP -> XRP
$ws.Range("E9").Value = "  -0.27%  "

# Row 10: Dogecoin -> Dogecoin
$ws.Range("E10").Value = "  +1.18%  "

# Row 11: Toncoin -> Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.30"
$ws.Range("E11").Value = "  -1.46%  "

# Row 12: Cardano -> Cardano
$ws.Range("E12").Value = "  -0.21%  "

# Row 13: ShibaInu -> ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").Value = "  -0.79%  "

# Row 14: Avalanche -> Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.19"
$ws.Range("E14").Value = "  +0.72%  "

# Row 15: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.446.97"
$ws.Range("E15").Value = "  +0.97%  "

# Row 16: WrappedEther -> WrappedEther
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.778.54"
$ws.Range("E16").Value = "  +1.34%  "

# Row 17: Chainlink -> Chainlink
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.60"
$ws.Range("E17").Value = "  +4.70%  "

# Row 18: WrappedBTC -> WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.804.28"
$ws.Range("E18").Value = "  +0.37%  "

# Row 19: Polkadot -> Polkadot
$ws.Range("E19").Value = "  +2.63%  "

# Row 20: TRON -> TRON
$ws.Range("E20").Value = "  +0.23%  "

# Row 21: BitcoinCash -> BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "461.84"
$ws.Range("E21").Value = "  +0.80%  "

# Row 22: Uniswap -> Uniswap
$ws.Range("E22").Value = "  -5.93%  "

# Row 23: Polygon -> Polygon
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.701"
$ws.Range("E23").Value = "  +0.69%  "

# Row 24: PEPE -> PEPE
$ws.Range("E24").Value = "  +0.34%  "

# Row 25: Litecoin -> Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.62"
$ws.Range("E25").Value = "  +0.36%  "

# Row 26: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.14"
$ws.Range("E26").Value = "  +2.57%  "

# Row 27: Fetch.AI -> Fetch.AI
$ws.Range("E27").Value = "  -1.44%  "

# Row 28: Dai -> RenderToken
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("E28").Value = "  -0.05%  "

# Row 29: RenderToken -> Dai
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.08%  "

# Row 30: WrappedeETH -> WrappedeETH
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.955.20"
$ws.Range("E30").Value = "  +0.94%  "

# Row 31: PancakeSwap -> PancakeSwap
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.80"
$ws.Range("E31").Value = "  +1.88%  "

# Row 32: ImmutableX -> ImmutableX
$ws.Range("E32").Value = "  +4.57%  "

# Row 33: NEARProtocol -> NEARProtocol
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.28"
$ws.Range("E33").Value = "  +1.03%  "

# Row 34: EthereumClassic -> EthereumClassic
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.73"
$ws.Range("E34").Value = "  +0.10%  "

# Row 35: Binance-PegBSC-USD -> Aptos
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.11"
$ws.Range("E35").Value = "  -0.21%  "

# Row 36: Aptos -> Binance-PegBSC-USD
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.10%  "

# Row 37: RenzoRestakedETH -> RenzoRestakedETH
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.748.84"
$ws.Range("E37").Value = "  +0.74%  "

# Row 38: Hedera -> Hedera
$ws.Range("E38").Value = "  +0.13%  "

# Row 39: dogwifhat -> dogwifhat
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.44"
$ws.Range("E39").Value = "  +0.86%  "

# Row 40: Kaspa -> Kaspa
$ws.Range("E40").Value = "  +0.45%  "

# Row 41: Mantle -> Mantle
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.54%  "

# Row 42: Filecoin -> Filecoin
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.80"
$ws.Range("E42").Value = "  +1.00%  "

# Row 43: FirstDigitalUSD -> FirstDigitalUSD
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.11%  "

# Row 45: OKB -> OKB
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.19"
$ws.Range("E45").Value = "  +2.63%  "

# Row 46: Arweave -> Arweave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.77"
$ws.Range("E46").Value = "  -1.05%  "

# Row 47: TheGraph -> TheGraph
$ws.Range("E47").Value = "  +0.04%  "

# Row 48: Monero -> Monero
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "149.47"
$ws.Range("E48").Value = "  +2.15%  "

# Row 49: Cosmos -> Cosmos
$ws.Range("E49").Value = "  +0.06%  "

# Row 50: Bittensor -> Bittensor
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "399.05"
$ws.Range("E50").Value = "  +2.29%  "

# Row 51: Stacks -> Stacks
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.84"
$ws.Range("E51").Value = "  -3.00%  "
